$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6
# from serial 45233 (2023-11-03) to 45243 (2023-11-13)
$ws.Range("C2:C6").Value = 45243
